$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row 7 of tracker data ---------------------------------
# Date value (Excel serial 42608 = 2016-08-26)
$ws.Cells.Item(7, 1).Value = 42608

# Columns B..J -> "Done"
$ws.Range("B7:J7").Value = "Done"

# Columns K..L -> "To do"
$ws.Range("K7:L7").Value = "To do"

# --- Copy formatting from existing rows so styles are re-used ----------
# A7 should look like the date cells above it (A4:A6)
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# B7:J7 should have the plain (no border/wrap) formatting used by B6:H6
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7:J7").PasteSpecial(-4122) | Out-Null

# K7:L7 should have the wrap-text formatting used by K6:L6
$ws.Range("K6:L6").Copy() | Out-Null
$ws.Range("K7:L7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Update the view/selection state ------------------------------------
$ws.Range("L7").Select() | Out-Null
